$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must remain a TEXT string even when it
# looks like a plain number (e.g. "23.63"), so Excel does not silently
# coerce it into a numeric cell (which would also introduce float noise
# such as 23.629999999999999).
function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '27.986.61'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').Value = '1.641.63'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('E7').Value = '  +0.02%  '
Set-TextValue 'D8' '23.63'
$ws.Range('E8').Value = '  +1.47%  '
$ws.Range('E9').Value = '  -1.07%  '
$ws.Range('E10').Value = '  +0.50%  '
$ws.Range('E11').Value = '  +2.25%  '
$ws.Range('D12').Value = '1.875.21'
$ws.Range('E12').Value = '  +0.52%  '
$ws.Range('D13').Value = '1.644.13'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D14' '0.576'
$ws.Range('E14').Value = '  +4.26%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D15' '4.10'
$ws.Range('E15').Value = '  +1.34%  '
Set-TextValue 'D16' '65.95'
$ws.Range('E16').Value = '  +1.18%  '
$ws.Range('D17').Value = '27.985.01'
$ws.Range('E17').Value = '  +1.30%  '
Set-TextValue 'D18' '232.33'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('D19').Value = '0.0₃0724'
$ws.Range('E19').Value = '  +0.54%  '
Set-TextValue 'D20' '7.62'
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('E21').Value = '  +0.02%  '
Set-TextValue 'D22' '10.77'
$ws.Range('E22').Value = '  +1.28%  '
Set-TextValue 'D23' '4.36'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E24').Value = '  -0.97%  '
$ws.Range('E25').Value = '  +1.89%  '
Set-TextValue 'D26' '6.97'
$ws.Range('E26').Value = '  +1.27%  '
Set-TextValue 'D27' '15.73'
$ws.Range('E27').Value = '  +1.41%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  +0.57%  '
$ws.Range('E31').Value = '  +0.19%  '
Set-TextValue 'D32' '3.35'
$ws.Range('E32').Value = '  +2.19%  '
$ws.Range('E33').Value = '  +0.58%  '
$ws.Range('D34').Value = '1.414.39'
$ws.Range('E34').Value = '  -4.28%  '
$ws.Range('E35').Value = '  +2.51%  '
$ws.Range('E36').Value = '  +0.94%  '
Set-TextValue 'D37' '0.888'
$ws.Range('E37').Value = '  +1.07%  '
$ws.Range('E38').Value = '  +0.98%  '
Set-TextValue 'D39' '0.558'
$ws.Range('E39').Value = '  -0.13%  '
Set-TextValue 'D40' '0.918'
$ws.Range('E40').Value = '  -4.02%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('E43').Value = '  +6.70%  '
Set-TextValue 'D44' '66.27'
$ws.Range('E44').Value = '  -2.17%  '
$ws.Range('E45').Value = '  +2.83%  '
$ws.Range('E46').Value = '  +0.18%  '
$ws.Range('D47').Value = '1.783.28'
$ws.Range('E47').Value = '  +0.53%  '
Set-TextValue 'D48' '88.07'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('E49').Value = '  +1.33%  '
Set-TextValue 'D51' '7.62'
$ws.Range('E51').Value = '  -1.28%  '
